$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = "b"
$ws.Range("J3").Value = "Acknowledge (Backchannel)"
$ws.Range("I5").Value = "sd"
$ws.Range("J5").Value = "Statement-non-opinion"
$ws.Range("I7").Value = "sd"
$ws.Range("J7").Value = "Statement-non-opinion"
$ws.Range("I12").Value = "aa"
$ws.Range("J12").Value = "Agree/Accept"
$ws.Range("I24").Value = "%"
$ws.Range("J24").Value = "Uninterpretable"
$ws.Range("I26").Value = "ba"
$ws.Range("J26").Value = "Appreciation"
$ws.Range("I27").Value = "sv"
$ws.Range("J27").Value = "Statement-opinion"
$ws.Range("I57").Value = "sv"
$ws.Range("J57").Value = "Statement-opinion"
$ws.Range("I58").Value = "ba"
$ws.Range("J58").Value = "Appreciation"
$ws.Range("I63").Value = "sv"
$ws.Range("J63").Value = "Statement-opinion"
$ws.Range("I65").Value = "sv"
$ws.Range("J65").Value = "Statement-opinion"
$ws.Range("I74").Value = "sd"
$ws.Range("J74").Value = "Statement-non-opinion"
$ws.Range("I77").Value = "sd"
$ws.Range("J77").Value = "Statement-non-opinion"
$ws.Range("I81").Value = "b"
$ws.Range("J81").Value = "Acknowledge (Backchannel)"
$ws.Range("I85").Value = "aa"
$ws.Range("J85").Value = "Agree/Accept"
$ws.Range("I87").Value = "aa"
$ws.Range("J87").Value = "Agree/Accept"
$ws.Range("I90").Value = "b"
$ws.Range("J90").Value = "Acknowledge (Backchannel)"
$ws.Range("I91").Value = "sv"
$ws.Range("J91").Value = "Statement-opinion"
$ws.Range("I97").Value = "sv"
$ws.Range("J97").Value = "Statement-opinion"
$ws.Range("I100").Value = "ba"
$ws.Range("J100").Value = "Appreciation"
$ws.Range("I111").Value = "sd"
$ws.Range("J111").Value = "Statement-non-opinion"
$ws.Range("I125").Value = "%"
$ws.Range("J125").Value = "Uninterpretable"
$ws.Range("I131").Value = "%"
$ws.Range("J131").Value = "Uninterpretable"
$ws.Range("I142").Value = "sv"
$ws.Range("J142").Value = "Statement-opinion"
$ws.Range("I143").Value = "sv"
$ws.Range("J143").Value = "Statement-opinion"
$ws.Range("I155").Value = "sd"
$ws.Range("J155").Value = "Statement-non-opinion"
$ws.Range("I160").Value = "sd"
$ws.Range("J160").Value = "Statement-non-opinion"
$ws.Range("I162").Value = "sd"
$ws.Range("J162").Value = "Statement-non-opinion"
$ws.Range("I173").Value = "%"
$ws.Range("J173").Value = "Uninterpretable"
$ws.Range("I175").Value = "%"
$ws.Range("J175").Value = "Uninterpretable"
$ws.Range("I194").Value = "sv"
$ws.Range("J194").Value = "Statement-opinion"
$ws.Range("I195").Value = "b"
$ws.Range("J195").Value = "Acknowledge (Backchannel)"
$ws.Range("I197").Value = "sd"
$ws.Range("J197").Value = "Statement-non-opinion"
$ws.Range("I198").Value = "sd"
$ws.Range("J198").Value = "Statement-non-opinion"
$ws.Range("I205").Value = "sd"
$ws.Range("J205").Value = "Statement-non-opinion"
$ws.Range("I208").Value = "aa"
$ws.Range("J208").Value = "Agree/Accept"
$ws.Range("I209").Value = "ba"
$ws.Range("J209").Value = "Appreciation"
$ws.Range("I215").Value = "aa"
$ws.Range("J215").Value = "Agree/Accept"
$ws.Range("I222").Value = "sd"
$ws.Range("J222").Value = "Statement-non-opinion"
$ws.Range("I230").Value = "aa"
$ws.Range("J230").Value = "Agree/Accept"
$ws.Range("I232").Value = "sd"
$ws.Range("J232").Value = "Statement-non-opinion"
$ws.Range("I235").Value = "%"
$ws.Range("J235").Value = "Uninterpretable"
$ws.Range("I247").Value = "sv"
$ws.Range("J247").Value = "Statement-opinion"
$ws.Range("I248").Value = "sd"
$ws.Range("J248").Value = "Statement-non-opinion"
$ws.Range("I256").Value = "sv"
$ws.Range("J256").Value = "Statement-opinion"
$ws.Range("I258").Value = "b"
$ws.Range("J258").Value = "Acknowledge (Backchannel)"
$ws.Range("I260").Value = "sd"
$ws.Range("J260").Value = "Statement-non-opinion"
$ws.Range("I261").Value = "sd"
$ws.Range("J261").Value = "Statement-non-opinion"
$ws.Range("I263").Value = "sd"
$ws.Range("J263").Value = "Statement-non-opinion"
$ws.Range("I265").Value = "sd"
$ws.Range("J265").Value = "Statement-non-opinion"
$ws.Range("I269").Value = "aa"
$ws.Range("J269").Value = "Agree/Accept"
$ws.Range("I274").Value = "b"
$ws.Range("J274").Value = "Acknowledge (Backchannel)"
$ws.Range("I278").Value = "sd"
$ws.Range("J278").Value = "Statement-non-opinion"
$ws.Range("I288").Value = "aa"
$ws.Range("J288").Value = "Agree/Accept"
$ws.Range("I302").Value = "sd"
$ws.Range("J302").Value = "Statement-non-opinion"
$ws.Range("I306").Value = "aa"
$ws.Range("J306").Value = "Agree/Accept"
$ws.Range("I317").Value = "sd"
$ws.Range("J317").Value = "Statement-non-opinion"
$ws.Range("I320").Value = "sv"
$ws.Range("J320").Value = "Statement-opinion"
$ws.Range("I322").Value = "ba"
$ws.Range("J322").Value = "Appreciation"
$ws.Range("I324").Value = "sv"
$ws.Range("J324").Value = "Statement-opinion"
$ws.Range("I352").Value = "sv"
$ws.Range("J352").Value = "Statement-opinion"
$ws.Range("I358").Value = "sv"
$ws.Range("J358").Value = "Statement-opinion"
$ws.Range("I363").Value = "aa"
$ws.Range("J363").Value = "Agree/Accept"
$ws.Range("I368").Value = "%"
$ws.Range("J368").Value = "Uninterpretable"
$ws.Range("I370").Value = "%"
$ws.Range("J370").Value = "Uninterpretable"
$ws.Range("I375").Value = "b"
$ws.Range("J375").Value = "Acknowledge (Backchannel)"
$ws.Range("I396").Value = "sv"
$ws.Range("J396").Value = "Statement-opinion"
$ws.Range("I397").Value = "b"
$ws.Range("J397").Value = "Acknowledge (Backchannel)"
$ws.Range("I404").Value = "ba"
$ws.Range("J404").Value = "Appreciation"
$ws.Range("I414").Value = "sv"
$ws.Range("J414").Value = "Statement-opinion"
$ws.Range("I416").Value = "ba"
$ws.Range("J416").Value = "Appreciation"
$ws.Range("I417").Value = "sd"
$ws.Range("J417").Value = "Statement-non-opinion"
$ws.Range("I422").Value = "ba"
$ws.Range("J422").Value = "Appreciation"
$ws.Range("I426").Value = "sd"
$ws.Range("J426").Value = "Statement-non-opinion"
$ws.Range("I428").Value = "%"
$ws.Range("J428").Value = "Uninterpretable"
$ws.Range("I445").Value = "ba"
$ws.Range("J445").Value = "Appreciation"
$ws.Range("I451").Value = "sv"
$ws.Range("J451").Value = "Statement-opinion"
$ws.Range("I452").Value = "sv"
$ws.Range("J452").Value = "Statement-opinion"
$ws.Range("I461").Value = "sd"
$ws.Range("J461").Value = "Statement-non-opinion"
$ws.Range("I471").Value = "sd"
$ws.Range("J471").Value = "Statement-non-opinion"
$ws.Range("I473").Value = "ba"
$ws.Range("J473").Value = "Appreciation"
$ws.Range("I477").Value = "sv"
$ws.Range("J477").Value = "Statement-opinion"
$ws.Range("I488").Value = "sd"
$ws.Range("J488").Value = "Statement-non-opinion"
$ws.Range("I491").Value = "b"
$ws.Range("J491").Value = "Acknowledge (Backchannel)"
$ws.Range("I494").Value = "%"
$ws.Range("J494").Value = "Uninterpretable"
$ws.Range("I496").Value = "aa"
$ws.Range("J496").Value = "Agree/Accept"
$ws.Range("I501").Value = "aa"
$ws.Range("J501").Value = "Agree/Accept"
$ws.Range("I503").Value = "aa"
$ws.Range("J503").Value = "Agree/Accept"
$ws.Range("I508").Value = "%"
$ws.Range("J508").Value = "Uninterpretable"
$ws.Range("I518").Value = "%"
$ws.Range("J518").Value = "Uninterpretable"
$ws.Range("I531").Value = "b"
$ws.Range("J531").Value = "Acknowledge (Backchannel)"
